# Digital_Signage_kayttajatestaus.xlsx edit script
# Commit message: "fixes header font and upload error images"
#
# Summary of changes (derived from the OOXML diff):
#  - "Myyrmäki" sheet:      D18 ("Sunnuntai = ei ruokalistaa") -> cleared;
#                           E18 status changed from "Ei suoritettu" (yellow) to "Hyväksytty" (green)
#  - "Karaportti" sheet:    D18 ("Sunnuntai = ei ruokalistaa") -> "Ruokalista ei näy";
#                           E18 status changed from "Ei suoritettu" (yellow) to "Hylätty" (red)
#  - "Myllypuro (ENG)":     D19 ("Sunnuntai = ei ruokalistaa") -> cleared;
#                           E19 status changed from "Ei suoritettu" (yellow) to "Hyväksytty" (green)
#  - "Arabia" sheet:        D18 ("Sunnuntai = ei ruokalistaa") -> "Oikeanlainen virheilmoitus näkyy";
#                           column D widened slightly
#  - Selection (active cell) changes on every sheet.

$wb = $excel.ActiveWorkbook

$wsMyyrmaki  = $wb.Worksheets.Item("Myyrmäki")
$wsKaraportti = $wb.Worksheets.Item("Karaportti")
$wsMyllypuro = $wb.Worksheets.Item("Myllypuro (ENG)")
$wsArabia    = $wb.Worksheets.Item("Arabia")

# ---------------------------------------------------------------------------
# Arabia: replace note with new error-message wording, widen column D
# (done first so the new shared-string entries land in the same order as
# the target workbook: "Oikeanlainen virheilmoitus näkyy" before
# "Ruokalista ei näy")
# ---------------------------------------------------------------------------
$wsArabia.Range("D18").Value2 = "Oikeanlainen virheilmoitus näkyy"
$wsArabia.Columns.Item(4).ColumnWidth = 28.333333333333332

# ---------------------------------------------------------------------------
# Karaportti: replace note with new error-image wording, mark row as failed
# ---------------------------------------------------------------------------
$wsKaraportti.Range("D18").Value2 = "Ruokalista ei näy"
$wsKaraportti.Range("E18").Value2 = "Hylätty"
$wsKaraportti.Range("D4").Copy()
$wsKaraportti.Range("E18").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# Myyrmäki: clear the "Sunnuntai = ei ruokalistaa" note, mark row as passed
# ---------------------------------------------------------------------------
$wsMyyrmaki.Range("D18").Value2 = ""
$wsMyyrmaki.Range("E18").Value2 = "Hyväksytty"
$wsMyyrmaki.Range("E14").Copy()
$wsMyyrmaki.Range("E18").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# Myllypuro (ENG): clear the "Sunnuntai = ei ruokalistaa" note, mark as passed
# ---------------------------------------------------------------------------
$wsMyllypuro.Range("D19").Value2 = ""
$wsMyllypuro.Range("E19").Value2 = "Hyväksytty"
$wsMyllypuro.Range("E14").Copy()
$wsMyllypuro.Range("E19").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# Restore / update the active-cell selection on every sheet. The workbook
# keeps "Myllypuro (ENG)" as the active tab, so it must be selected last.
# ---------------------------------------------------------------------------
$wsMyyrmaki.Range("D18").Select()
$wsKaraportti.Range("F30").Select()
$wsArabia.Range("H18").Select()
$wsMyllypuro.Range("G16").Select()
